$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("U2").Value = 2.02
$ws.Range("AA3").Value = 980
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 980
$ws.Range("AF3").Value = 980
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 980
$ws.Range("AJ3").Value = 980
$ws.Range("AK3").Value = 980
$ws.Range("AL3").Value = 980
$ws.Range("AM3").Value = 75
$ws.Range("AN3").Value = 980
$ws.Range("AO3").Value = 1000
$ws.Range("G3").Value = 2.84
$ws.Range("I3").Value = 2.76
$ws.Range("J3").Value = 3.65
$ws.Range("L3").Value = 1.28
$ws.Range("N3").Value = 4.5
$ws.Range("P3").Value = 2.24
$ws.Range("R3").Value = 1.49
$ws.Range("U3").Value = 2.24
$ws.Range("W3").Value = 1.54
$ws.Range("X3").Value = 980
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 980
$ws.Range("N5").Value = 3.1
$ws.Range("P5").Value = 1.71
$ws.Range("S5").Value = 3.7
$ws.Range("AB6").Value = 6.8
$ws.Range("AC6").Value = 1000
$ws.Range("AH6").Value = 36
$ws.Range("AN6").Value = 1000
$ws.Range("U6").Value = 1.54
$ws.Range("V6").Value = 1.09
$ws.Range("X6").Value = 1000
$ws.Range("AD7").Value = 1000
$ws.Range("F7").Value = 1.92
$ws.Range("U7").Value = 1.8
$ws.Range("I8").Value = 8.4
$ws.Range("G9").Value = 2.52
$ws.Range("J9").Value = 3.4
$ws.Range("V9").Value = 1.35
$ws.Range("G10").Value = 2.72
$ws.Range("U10").Value = 1.94
$ws.Range("N11").Value = 3
$ws.Range("O11").Value = 1.4
$ws.Range("R11").Value = 1.22
$ws.Range("AB12").Value = 10.5
$ws.Range("AC12").Value = 10
$ws.Range("AD12").Value = 23
$ws.Range("AE12").Value = 80
$ws.Range("AF12").Value = 16.5
$ws.Range("AG12").Value = 14.5
$ws.Range("AH12").Value = 26
$ws.Range("AI12").Value = 90
$ws.Range("AJ12").Value = 36
$ws.Range("AK12").Value = 34
$ws.Range("F12").Value = 2.02
$ws.Range("G12").Value = 2.28
$ws.Range("H12").Value = 3.8
$ws.Range("I12").Value = 4.7
$ws.Range("J12").Value = 3.15
$ws.Range("K12").Value = 4
$ws.Range("L12").Value = 1.35
$ws.Range("M12").Value = 1.08
$ws.Range("N12").Value = 3.05
$ws.Range("O12").Value = 1.39
$ws.Range("P12").Value = 1.68
$ws.Range("Q12").Value = 2.04
$ws.Range("R12").Value = 1.23
$ws.Range("S12").Value = 3.45
$ws.Range("T12").Value = 1.79
$ws.Range("U12").Value = 1.82
$ws.Range("W12").Value = 1.78
$ws.Range("X12").Value = 1000
$ws.Range("Y12").Value = 17
$ws.Range("Z12").Value = 38
$ws.Range("AA13").Value = 90
$ws.Range("AB13").Value = 11.5
$ws.Range("AC13").Value = 9.800000000000001
$ws.Range("AD13").Value = 20
$ws.Range("AE13").Value = 60
$ws.Range("AF13").Value = 20
$ws.Range("AG13").Value = 15.5
$ws.Range("AH13").Value = 26
$ws.Range("AI13").Value = 85
$ws.Range("AJ13").Value = 48
$ws.Range("AK13").Value = 42
$ws.Range("AL13").Value = 70
$ws.Range("F13").Value = 2.22
$ws.Range("G13").Value = 2.9
$ws.Range("I13").Value = 4.3
$ws.Range("K13").Value = 4.3
$ws.Range("R13").Value = 1.19
$ws.Range("S13").Value = 3.75
$ws.Range("V13").Value = 1.3
$ws.Range("X13").Value = 14.5
$ws.Range("Y13").Value = 15
$ws.Range("Z13").Value = 30
$ws.Range("H15").Value = 3.8
$ws.Range("I16").Value = 4.2
$ws.Range("J16").Value = 4.3
$ws.Range("Q16").Value = 1.31
$ws.Range("R16").Value = 1.83
$ws.Range("V16").Value = 1.31
$ws.Range("AA17").Value = 27
$ws.Range("AB17").Value = 25
$ws.Range("AC17").Value = 12
$ws.Range("AD17").Value = 13
$ws.Range("AE17").Value = 21
$ws.Range("AF17").Value = 38
$ws.Range("AG17").Value = 19
$ws.Range("AH17").Value = 18
$ws.Range("AI17").Value = 32
$ws.Range("AJ17").Value = 80
$ws.Range("AK17").Value = 46
$ws.Range("AL17").Value = 46
$ws.Range("AM17").Value = 70
$ws.Range("AN17").Value = 32
$ws.Range("AO17").Value = 10
$ws.Range("N17").Value = 5.4
$ws.Range("O17").Value = 1.18
$ws.Range("R17").Value = 1.61
$ws.Range("S17").Value = 2.3
$ws.Range("T17").Value = 1.54
$ws.Range("U17").Value = 2.34
$ws.Range("X17").Value = 30
$ws.Range("Y17").Value = 16
$ws.Range("Z17").Value = 18
$ws.Range("P18").Value = 1.58
$ws.Range("F19").Value = 2.48
$ws.Range("F21").Value = 2.94
$ws.Range("H21").Value = 2.96
$ws.Range("I21").Value = 2.98
$ws.Range("J21").Value = 3.1
$ws.Range("K21").Value = 3.15
$ws.Range("N21").Value = 3
$ws.Range("O21").Value = 1.48
$ws.Range("V21").Value = 1.5
$ws.Range("Y21").Value = 9.199999999999999
$ws.Range("F22").Value = 1.94
$ws.Range("G22").Value = 1.95
$ws.Range("L22").Value = 1.28
$ws.Range("N22").Value = 6.2
$ws.Range("P22").Value = 2.82
$ws.Range("Q22").Value = 1.53
$ws.Range("T22").Value = 1.54
$ws.Range("AJ23").Value = 32
$ws.Range("F23").Value = 2.34
$ws.Range("G23").Value = 2.36
$ws.Range("I23").Value = 3.55
$ws.Range("L23").Value = 1.48
$ws.Range("S23").Value = 4.3
$ws.Range("V23").Value = 1.39
$ws.Range("AM24").Value = 90
$ws.Range("AN24").Value = 29
$ws.Range("G24").Value = 3.05
$ws.Range("H24").Value = 2.58
$ws.Range("W24").Value = 1.48
$ws.Range("AC25").Value = 8.800000000000001
$ws.Range("AL25").Value = 65
$ws.Range("AN25").Value = 980
$ws.Range("AO25").Value = 80
$ws.Range("L25").Value = 1.48
$ws.Range("Z25").Value = 980
